$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45175 = 2023-09-06) for every
# data row (2..210). The update bumps that date to serial 45177 (2023-09-08)
# for all of those rows.
$ws.Range("C2:C210").Value = 45177
